# Daily scrape update - 2025-10-25 03:15:35 UTC
# Updates existing listing rows 2-7 with refreshed data, appends new
# listing rows 8-20, and widens columns C, F (narrower), H to fit new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (OPPORTUNITY ID) holds numeric-looking values that must stay
# text, exactly like the existing header rows. Pre-format as Text so
# Excel does not silently convert the IDs to numbers.
$ws.Range("A2:A20").NumberFormat = "@"

# --- Refresh existing rows 2-7 (only the cells whose values changed) ---
$ws.Cells.Item(2, 1).Value = '1328782'
$ws.Cells.Item(2, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328782'
$ws.Cells.Item(2, 3).Value = 'Risk Intern'
$ws.Cells.Item(2, 6).Value = '0 applicants'
$ws.Cells.Item(2, 8).Value = 'Banesco Panamá'
$ws.Cells.Item(3, 1).Value = '1328781'
$ws.Cells.Item(3, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328781'
$ws.Cells.Item(3, 3).Value = 'Risk Trainee'
$ws.Cells.Item(3, 8).Value = 'Banesco Panamá'
$ws.Cells.Item(4, 1).Value = '1328780'
$ws.Cells.Item(4, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328780'
$ws.Cells.Item(4, 3).Value = 'HR Transformation Intern'
$ws.Cells.Item(4, 4).Value = 'Panamá, Provincia de Panamá, Panamá'
$ws.Cells.Item(4, 6).Value = '0 applicants'
$ws.Cells.Item(4, 7).Value = '6 - 18 Months'
$ws.Cells.Item(4, 8).Value = 'HILTI Panama'
$ws.Cells.Item(5, 1).Value = '1328779'
$ws.Cells.Item(5, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328779'
$ws.Cells.Item(5, 3).Value = 'AMS Regional Support Intern'
$ws.Cells.Item(5, 4).Value = 'Panamá, Provincia de Panamá, Panamá'
$ws.Cells.Item(5, 6).Value = '0 applicants'
$ws.Cells.Item(5, 7).Value = '6 - 18 Months'
$ws.Cells.Item(5, 8).Value = 'HILTI Panama'
$ws.Cells.Item(6, 1).Value = '1328774'
$ws.Cells.Item(6, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328774'
$ws.Cells.Item(6, 3).Value = 'Digital Marketing'
$ws.Cells.Item(6, 4).Value = 'Leiria, Portugal'
$ws.Cells.Item(6, 6).Value = '1 applicant'
$ws.Cells.Item(6, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(6, 8).Value = 'Multidrive'
$ws.Cells.Item(7, 1).Value = '1328773'
$ws.Cells.Item(7, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328773'
$ws.Cells.Item(7, 3).Value = 'Procurement Management Assistant Intern'
$ws.Cells.Item(7, 4).Value = 'Panamá, Provincia de Panamá, Panamá'
$ws.Cells.Item(7, 6).Value = '3 applicants'
$ws.Cells.Item(7, 8).Value = 'Samsung Electronics Latinoamerica SELA'

# --- Append new rows 8-20 scraped today ---
# Row 8
$ws.Cells.Item(8, 1).Value = '1328768'
$ws.Cells.Item(8, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328768'
$ws.Cells.Item(8, 3).Value = 'Accelerate Romania - Digital Content Intern – Social Media & Website'
$ws.Cells.Item(8, 4).Value = 'Cluj-Napoca, Romania'
$ws.Cells.Item(8, 5).Value = 'No'
$ws.Cells.Item(8, 6).Value = '3 applicants'
$ws.Cells.Item(8, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(8, 8).Value = 'Dog Assist'

# Row 9
$ws.Cells.Item(9, 1).Value = '1328767'
$ws.Cells.Item(9, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328767'
$ws.Cells.Item(9, 3).Value = 'Accelerate Romania - Digital Marketing & Social Media Assistant'
$ws.Cells.Item(9, 4).Value = 'Cluj-Napoca, Romania'
$ws.Cells.Item(9, 5).Value = 'No'
$ws.Cells.Item(9, 6).Value = '0 applicants'
$ws.Cells.Item(9, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(9, 8).Value = 'Ave Visto'

# Row 10
$ws.Cells.Item(10, 1).Value = '1328766'
$ws.Cells.Item(10, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328766'
$ws.Cells.Item(10, 3).Value = 'Accelerate Romania - Community Manager'
$ws.Cells.Item(10, 4).Value = 'Cluj-Napoca, Romania'
$ws.Cells.Item(10, 5).Value = 'No'
$ws.Cells.Item(10, 6).Value = '0 applicants'
$ws.Cells.Item(10, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(10, 8).Value = 'ClujStartups'

# Row 11
$ws.Cells.Item(11, 1).Value = '1328765'
$ws.Cells.Item(11, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328765'
$ws.Cells.Item(11, 3).Value = 'Accelerate Romania -Business Development Specialist'
$ws.Cells.Item(11, 4).Value = 'Cluj-Napoca, Romania'
$ws.Cells.Item(11, 5).Value = 'No'
$ws.Cells.Item(11, 6).Value = '1 applicant'
$ws.Cells.Item(11, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(11, 8).Value = 'MTD Technology'

# Row 12
$ws.Cells.Item(12, 1).Value = '1328764'
$ws.Cells.Item(12, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328764'
$ws.Cells.Item(12, 3).Value = 'Accelerate Romania - Market Development Intern'
$ws.Cells.Item(12, 4).Value = 'Brașov, Romania'
$ws.Cells.Item(12, 5).Value = 'No'
$ws.Cells.Item(12, 6).Value = '0 applicants'
$ws.Cells.Item(12, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(12, 8).Value = 'Biklo'

# Row 13
$ws.Cells.Item(13, 1).Value = '1328763'
$ws.Cells.Item(13, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328763'
$ws.Cells.Item(13, 3).Value = 'Accelerate Romania - Social Media Manager'
$ws.Cells.Item(13, 4).Value = 'Cluj-Napoca, Romania'
$ws.Cells.Item(13, 5).Value = 'No'
$ws.Cells.Item(13, 6).Value = '0 applicants'
$ws.Cells.Item(13, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(13, 8).Value = 'nclav'

# Row 14
$ws.Cells.Item(14, 1).Value = '1328762'
$ws.Cells.Item(14, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328762'
$ws.Cells.Item(14, 3).Value = 'Accelerate - Marketing Specialist'
$ws.Cells.Item(14, 4).Value = 'Cluj-Napoca, Romania'
$ws.Cells.Item(14, 5).Value = 'No'
$ws.Cells.Item(14, 6).Value = '1 applicant'
$ws.Cells.Item(14, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(14, 8).Value = 'Emotionstudios'

# Row 15
$ws.Cells.Item(15, 1).Value = '1328761'
$ws.Cells.Item(15, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328761'
$ws.Cells.Item(15, 3).Value = 'Accelerate Romania - Online Marketing Intern'
$ws.Cells.Item(15, 4).Value = 'Sibiu, Romania'
$ws.Cells.Item(15, 5).Value = 'No'
$ws.Cells.Item(15, 6).Value = '1 applicant'
$ws.Cells.Item(15, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(15, 8).Value = 'Gads Online Marketing'

# Row 16
$ws.Cells.Item(16, 1).Value = '1328760'
$ws.Cells.Item(16, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328760'
$ws.Cells.Item(16, 3).Value = 'Accelerate Romania - Sales Business Manager'
$ws.Cells.Item(16, 4).Value = 'Sibiu, Romania'
$ws.Cells.Item(16, 5).Value = 'No'
$ws.Cells.Item(16, 6).Value = '0 applicants'
$ws.Cells.Item(16, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(16, 8).Value = 'Gads Online Marketing'

# Row 17
$ws.Cells.Item(17, 1).Value = '1328759'
$ws.Cells.Item(17, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328759'
$ws.Cells.Item(17, 3).Value = 'Accelerate Romania - Graphic Design'
$ws.Cells.Item(17, 4).Value = 'Cluj-Napoca, Romania'
$ws.Cells.Item(17, 5).Value = 'No'
$ws.Cells.Item(17, 6).Value = '0 applicants'
$ws.Cells.Item(17, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(17, 8).Value = 'ROTSA'

# Row 18
$ws.Cells.Item(18, 1).Value = '1328758'
$ws.Cells.Item(18, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328758'
$ws.Cells.Item(18, 3).Value = 'Accelerate Romania - Business Development Associate'
$ws.Cells.Item(18, 4).Value = 'Cluj-Napoca, Romania'
$ws.Cells.Item(18, 5).Value = 'No'
$ws.Cells.Item(18, 6).Value = '0 applicants'
$ws.Cells.Item(18, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(18, 8).Value = 'Mejix'

# Row 19
$ws.Cells.Item(19, 1).Value = '1328750'
$ws.Cells.Item(19, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328750'
$ws.Cells.Item(19, 3).Value = 'Accelerate Serbia | Digital Marketing & Business Trainee'
$ws.Cells.Item(19, 4).Value = 'Београд, Србија'
$ws.Cells.Item(19, 5).Value = 'No'
$ws.Cells.Item(19, 6).Value = '1 applicant'
$ws.Cells.Item(19, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(19, 8).Value = 'TERMOPLUS d.o.o.'

# Row 20
$ws.Cells.Item(20, 1).Value = '1321400'
$ws.Cells.Item(20, 2).Value = 'https://aiesec.org/opportunity/global-talent/1321400'
$ws.Cells.Item(20, 3).Value = 'Digital Marketing Executive'
$ws.Cells.Item(20, 4).Value = 'Cairo, Cairo Governorate, Egypt'
$ws.Cells.Item(20, 5).Value = 'No'
$ws.Cells.Item(20, 6).Value = '33 applicants'
$ws.Cells.Item(20, 7).Value = '3 - 6 Months'
$ws.Cells.Item(20, 8).Value = 'Silverkey Technologies Egypt'

# --- Column width adjustments ---
$ws.Columns.Item(3).ColumnWidth = 70.16666666666667   # C: 67 -> 71
$ws.Columns.Item(6).ColumnWidth = 15.166666666666666   # F: 17 -> 16
$ws.Columns.Item(8).ColumnWidth = 40.166666666666664   # H: 35 -> 41

